$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 199
$ws.Range("B3").Value = 179
$ws.Range("B4").Value = 157
$ws.Range("B5").Value = 157
$ws.Range("B6").Value = 157
$ws.Range("B7").Value = 144
$ws.Range("B8").Value = 137
$ws.Range("B9").Value = 134
$ws.Range("B10").Value = 132
$ws.Range("B11").Value = 129
